$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text/string updates (safe to assign directly) ---
$ws.Range('D2').Value = '28.028.06'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '1.871.22'
$ws.Range('E3').Value = '  -0.18%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('E7').Value = '  +2.07%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('E9').Value = '  -3.61%  '
$ws.Range('E10').Value = '  -0.30%  '
$ws.Range('B11').Value = 'Polkadot'
$ws.Range('C11').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('E11').Value = '  -1.78%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('E12').Value = '  -0.48%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.867.69'
$ws.Range('E13').Value = '  -0.37%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('E14').Value = '  +1.67%  '
$ws.Range('B15').Value = 'BinanceUSD'
$ws.Range('C15').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('E15').Value = '  +0.08%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('E17').Value = '  -0.10%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('E19').Value = '  -1.88%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('E21').Value = '  -1.06%  '
$ws.Range('B22').Value = 'WrappedBTC'
$ws.Range('C22').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D22').Value = '28.069.27'
$ws.Range('E22').Value = '  -0.25%  '
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('E23').Value = '  -2.65%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('B25').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C25').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D25').Value = '2.079.36'
$ws.Range('E25').Value = '  -0.53%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('E26').Value = '  -2.65%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('E27').Value = '  +0.32%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('E28').Value = '  -0.96%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('E29').Value = '  -1.11%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('E30').Value = '  +0.89%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('E31').Value = '  -2.70%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('E32').Value = '  +2.91%  '
$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('B34').Value = 'FraxShare'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E34').Value = '  -1.62%  '
$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('E35').Value = '  -0.95%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('E36').Value = '  -0.83%  '
$ws.Range('B37').Value = 'Algorand'
$ws.Range('C37').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('E37').Value = '  +0.92%  '
$ws.Range('B38').Value = 'TheSandbox'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('E38').Value = '  +3.30%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('E39').Value = '  -0.64%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('E40').Value = '  +2.56%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('E41').Value = '  -2.41%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('E42').Value = '  -2.50%  '
$ws.Range('B43').Value = 'Decentraland'
$ws.Range('C43').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('E43').Value = '  +2.31%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E44').Value = '  -0.31%  '
$ws.Range('B45').Value = 'WEMIXTOKEN'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('E47').Value = '  +1.91%  '
$ws.Range('B48').Value = 'EOS'
$ws.Range('C48').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('E48').Value = '  -0.50%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('E49').Value = '  -0.37%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('E50').Value = '  -2.38%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('E51').Value = '  -1.21%  '

# --- Numeric-looking text updates (must force Text format to avoid Excel auto-converting to a number) ---
$numericFixCells = @{
    'D5' = '312.24'
    'D6' = '1.003'
    'D7' = '0.5156'
    'D8' = '0.3845'
    'D11' = '6.205'
    'D12' = '20.56'
    'D14' = '7.315'
    'D15' = '1.004'
    'D16' = '0.00001099'
    'D17' = '90.87'
    'D18' = '0.06641'
    'D19' = '17.74'
    'D20' = '1.003'
    'D21' = '6.033'
    'D23' = '11.10'
    'D24' = '2.256'
    'D26' = '2.509'
    'D27' = '157.57'
    'D28' = '20.54'
    'D29' = '124.97'
    'D30' = '0.1065'
    'D31' = '1.034'
    'D32' = '5.779'
    'D33' = '3.596'
    'D34' = '9.461'
    'D35' = '0.02424'
    'D36' = '0.06527'
    'D37' = '0.2200'
    'D38' = '0.6596'
    'D39' = '1.203'
    'D40' = '5.023'
    'D41' = '1.211'
    'D42' = '11.22'
    'D43' = '0.6143'
    'D44' = '13.13'
    'D45' = '1.283'
    'D46' = '3.678'
    'D47' = '2.028'
    'D48' = '1.219'
    'D49' = '121.12'
    'D50' = '78.60'
    'D51' = '0.06831'
}

foreach ($cellAddr in $numericFixCells.Keys) {
    $c = $ws.Range($cellAddr)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $numericFixCells[$cellAddr]
    $c.Style = $origStyle
}
